$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date (column C) from 2023-09-14 (45183) to 2023-09-15 (45184)
# for rows 2 through 7.
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}
